# Add a "Save" column (H) to the s_vals sheet, matching the existing
# header styling (bold, bordered, centered) used by the other headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the header style already applied to G1 by copying its formatting
# onto H1 before writing the new header text (keeps the same style index
# instead of minting a new, duplicate one).
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data values for the "Save" column.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
